$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.336.33"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.711.13"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.40"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5288"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06688"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2664"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.87"
$ws.Range("E10").Value = "  -4.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07670"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.504"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").Value = "1.946.91"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "1.715.01"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5831"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "0.0₅8218"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.05"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "27.362.06"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "222.76"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.628"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.38"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.017"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.84"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.686"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1205"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.237"
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.25"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05352"
$ws.Range("E30").Value = "  -4.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.288"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.466"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.436"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.876"
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9510"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.397"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5849"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01631"
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "1.103.27"
$ws.Range("E40").Value = "  +4.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.795"
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8397"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.92"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "1.854.10"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.70"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4538"
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.094"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05235"
$ws.Range("E51").Value = "  -0.38%  "
